# "Remove vendor" action (now guarded by a confirmation message box in the
# app UI before this runs) on the "vendors" sheet: the row for
# "Test vendor" (A2) is removed from the list. The remaining vendor names
# shift up by one row and the now-vacated last row is cleared, matching how
# the app's "Remove" button re-renders the list after a confirmed deletion.
# (The confirmation dialog itself is UI-only and has no effect on the saved
# workbook, so it isn't reproduced here in this headless run.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

$vendorRange = $ws.Range("A2:A5")
$values = @($vendorRange.Value())

# Shift each remaining vendor name up by one row, dropping the removed
# "Test vendor" entry, and clear the vacated final row.
for ($i = 0; $i -lt ($values.Length - 1); $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i + 1]
}
$ws.Cells.Item(5, 1).ClearContents()
